# Auto update Excel log
# Appends newly-logged sensor readings (2026-01-28 afternoon batch) to the
# PIR, Humidity and Temperature sheets of the SeniorConnect master log.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param(
        $Sheet,
        $Rows
    )
    foreach ($r in $Rows) {
        $rowNum = $r[0]
        $Sheet.Cells.Item($rowNum, 1).Value = "'" + $r[1]   # force text, not a date serial
        $Sheet.Cells.Item($rowNum, 2).Value = $r[2]
        $Sheet.Cells.Item($rowNum, 3).Value = $r[3]
        $Sheet.Cells.Item($rowNum, 4).Value = $r[4]
        $val = $r[5]
        if ($val -like "*%") { $val = "'" + $val }       # force text, not a percentage number
        $Sheet.Cells.Item($rowNum, 5).Value = $val
        $Sheet.Cells.Item($rowNum, 6).Value = $r[6]
    }
}

$pirRows = @(
    @(40,"2026-01-28","16:13:44","16:00","Bathroom","No Motion","Inactive"),
    @(41,"2026-01-28","16:13:45","16:00","Bathroom","No Motion","Inactive"),
    @(42,"2026-01-28","16:13:48","16:00","Bathroom","No Motion","Inactive"),
    @(43,"2026-01-28","16:13:53","16:00","Bathroom","No Motion","Inactive"),
    @(44,"2026-01-28","16:13:58","16:00","Bathroom","No Motion","Inactive"),
    @(45,"2026-01-28","16:14:03","16:00","Bathroom","No Motion","Inactive"),
    @(46,"2026-01-28","16:14:08","16:00","Bathroom","No Motion","Inactive"),
    @(47,"2026-01-28","16:14:14","16:00","Bathroom","No Motion","Inactive"),
    @(48,"2026-01-28","16:14:19","16:00","Bathroom","No Motion","Inactive"),
    @(49,"2026-01-28","16:14:24","16:00","Bathroom","No Motion","Inactive"),
    @(50,"2026-01-28","16:14:29","16:00","Bathroom","No Motion","Inactive"),
    @(51,"2026-01-28","16:14:34","16:00","Bathroom","No Motion","Inactive"),
    @(52,"2026-01-28","16:14:39","16:00","Bathroom","No Motion","Inactive"),
    @(53,"2026-01-28","16:14:44","16:00","Bathroom","No Motion","Inactive")
)

$humidityRows = @(
    @(40,"2026-01-28","16:13:45","16:00","Bathroom","88.3%","Active"),
    @(41,"2026-01-28","16:13:46","16:00","Bathroom","88.3%","Active"),
    @(42,"2026-01-28","16:13:50","16:00","Bathroom","87.3%","Active"),
    @(43,"2026-01-28","16:13:54","16:00","Bathroom","88.3%","Active"),
    @(44,"2026-01-28","16:13:58","16:00","Bathroom","87.3%","Active"),
    @(45,"2026-01-28","16:14:06","16:00","Bathroom","88.2%","Active"),
    @(46,"2026-01-28","16:14:10","16:00","Bathroom","87.3%","Active"),
    @(47,"2026-01-28","16:14:14","16:00","Bathroom","88.3%","Active"),
    @(48,"2026-01-28","16:14:22","16:00","Bathroom","87.4%","Active"),
    @(49,"2026-01-28","16:14:30","16:00","Bathroom","87.4%","Active"),
    @(50,"2026-01-28","16:14:34","16:00","Bathroom","88.3%","Active"),
    @(51,"2026-01-28","16:14:42","16:00","Bathroom","88.3%","Active")
)

$temperatureRows = @(
    @(40,"2026-01-28","16:13:45","16:00","Bathroom","22.8C","Active"),
    @(41,"2026-01-28","16:13:46","16:00","Bathroom","22.8C","Active"),
    @(42,"2026-01-28","16:13:50","16:00","Bathroom","22.8C","Active"),
    @(43,"2026-01-28","16:13:54","16:00","Bathroom","22.8C","Active"),
    @(44,"2026-01-28","16:13:58","16:00","Bathroom","22.8C","Active"),
    @(45,"2026-01-28","16:14:06","16:00","Bathroom","22.8C","Active"),
    @(46,"2026-01-28","16:14:10","16:00","Bathroom","22.8C","Active"),
    @(47,"2026-01-28","16:14:14","16:00","Bathroom","22.8C","Active"),
    @(48,"2026-01-28","16:14:22","16:00","Bathroom","22.8C","Active"),
    @(49,"2026-01-28","16:14:30","16:00","Bathroom","22.8C","Active"),
    @(50,"2026-01-28","16:14:34","16:00","Bathroom","22.8C","Active"),
    @(51,"2026-01-28","16:14:42","16:00","Bathroom","22.8C","Active")
)

$wsPir = $wb.Worksheets.Item("PIR")
Add-LogRows $wsPir $pirRows

$wsHumidity = $wb.Worksheets.Item("Humidity")
Add-LogRows $wsHumidity $humidityRows

$wsTemperature = $wb.Worksheets.Item("Temperature")
Add-LogRows $wsTemperature $temperatureRows

